# "Redesigend UI, fixed some minor Live Editor bugs"
# The ShopText sheet's data table (A1:F35, header in row 1) is re-sorted by
# column C ("height") descending. This reshuffles rows 2-27 (rows 2-35 that
# have no height value sort to the bottom); rows 28-35 are already blank in
# column C so their relative order/content is unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Perform the actual re-sort of the data (this is what moves/rewrites the
# cell values for rows 2-27 to match the new column-C-descending order).
$dataRange = $ws.Range("A1:F35")
$sortKey   = $ws.Range("C2:C35")
$dataRange.Sort($sortKey, 2, $null, $null, 1, $null, 1, 1)

# Record the sort as a persisted Sort (Data > Sort) configuration, the way
# using the Sort dialog/ribbon (rather than a one-off Range.Sort) leaves a
# <sortState> behind on the worksheet.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("C2:C35"), 0, 2, $null, 0)
$sortObj.SetRange($ws.Range("A1:F35"))
$sortObj.Header = 1
$sortObj.Apply()

# UI bits: active cell moved to J19, and the page orientation was touched
# (set to Portrait) in the page setup dialog.
$ws.Range("J19").Select() | Out-Null
$ws.PageSetup.Orientation = 1
